$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for rows 8-11 (subjects 6-9), columns B (Condition) through Q (Q15)
$data = @(
    @(8,  "G", 3, 3, 2, 1, 2, 2, 4, 4, 3, 3, 3, 3, 2, 3, 5),
    @(9,  "F", 1, 3, 3, 1, 2, 1, 4, 4, 2, 3, 2, 1, 2, 3, 6),
    @(10, "G", 5, 4, 5, 5, 6, 4, 5, 4, 4, 3, 5, 5, 5, 3, 6),
    @(11, "F", 4, 1, 3, 2, 3, 3, 5, 5, 6, 4, 5, 4, 6, 4, 6)
)

foreach ($rowData in $data) {
    $row = $rowData[0]
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

# Update the sheet view: remove topLeftCell freeze and change selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("R11").Select()
